# Generate Report for Handback
# Update the timestamp values recorded on the handback status report.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date (row for
# f1fd48bb-ecd5-45a8-aa2b-2572d2da1df2.md)
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-29 20:59:10"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for f1fd48bb-ecd5-45a8-aa2b-2572d2da1df2.*.zh-cn.xlf
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-29 20:58:59"
$wsZhCn.Range("K4").Value = "2016-08-29 20:59:44"

# "de-de" sheet: Correspond Handoff Datetime (same value as the Overview's
# Latest HO Xliff Generate Date) / Correspond Handback DateTime for
# f1fd48bb-ecd5-45a8-aa2b-2572d2da1df2.*.de-de.xlf
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-29 20:59:10"
$wsDeDe.Range("K4").Value = "2016-08-29 20:59:51"
